$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B3 value: add the leading quotation mark that was missing
$ws.Range("B3").Value = '"So do not fear, for I am with you..."'

# Update the active selection to B3 as shown in the diff
$ws.Range("B3").Select()
